$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap F:V between row 63 and row 65
$ws.Range("F63").Value = "Ath Bilbao"
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = "Getafe"
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 1.71
$ws.Range("K63").Value = "17/09/2023 09:02"
$ws.Range("L63").Value = 1.53
$ws.Range("M63").Value = "27/09/2023 18:31"
$ws.Range("N63").Value = 3.42
$ws.Range("O63").Value = "17/09/2023 09:02"
$ws.Range("P63").Value = 4.06
$ws.Range("Q63").Value = "27/09/2023 18:49"
$ws.Range("R63").Value = 5.44
$ws.Range("S63").Value = "17/09/2023 09:02"
$ws.Range("T63").Value = 7.73
$ws.Range("U63").Value = "27/09/2023 18:49"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/"
$ws.Range("F65").Value = "Villarreal"
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = "Girona"
$ws.Range("I65").Value = 2
$ws.Range("J65").Value = 1.71
$ws.Range("K65").Value = "17/09/2023 09:02"
$ws.Range("L65").Value = 2.17
$ws.Range("M65").Value = "27/09/2023 18:51"
$ws.Range("N65").Value = 4.14
$ws.Range("O65").Value = "17/09/2023 09:02"
$ws.Range("P65").Value = 3.72
$ws.Range("Q65").Value = "27/09/2023 18:51"
$ws.Range("R65").Value = 4.72
$ws.Range("S65").Value = "17/09/2023 09:02"
$ws.Range("T65").Value = 3.42
$ws.Range("U65").Value = "27/09/2023 18:51"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/spain/laliga/villarreal-girona/80EuTg3A/"

# Swap F:V between row 66 and row 67
$ws.Range("F66").Value = "Cadiz CF"
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = "Rayo Vallecano"
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 2.59
$ws.Range("K66").Value = "17/09/2023 09:02"
$ws.Range("L66").Value = 2.72
$ws.Range("M66").Value = "27/09/2023 21:19"
$ws.Range("N66").Value = 3.06
$ws.Range("O66").Value = "17/09/2023 09:02"
$ws.Range("P66").Value = 3.11
$ws.Range("Q66").Value = "27/09/2023 21:17"
$ws.Range("R66").Value = 3.11
$ws.Range("S66").Value = "17/09/2023 09:02"
$ws.Range("T66").Value = 3
$ws.Range("U66").Value = "27/09/2023 21:30"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/spain/laliga/cadiz-rayo-vallecano/CEYt8hRp/"
$ws.Range("F67").Value = "Valencia"
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = "Real Sociedad"
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = 2.33
$ws.Range("K67").Value = "17/09/2023 09:02"
$ws.Range("L67").Value = 2.57
$ws.Range("M67").Value = "27/09/2023 21:27"
$ws.Range("N67").Value = 3.14
$ws.Range("O67").Value = "17/09/2023 09:02"
$ws.Range("P67").Value = 3.02
$ws.Range("Q67").Value = "27/09/2023 21:27"
$ws.Range("R67").Value = 3.26
$ws.Range("S67").Value = "17/09/2023 09:02"
$ws.Range("T67").Value = 3.31
$ws.Range("U67").Value = "27/09/2023 21:27"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/spain/laliga/valencia-real-sociedad/M3IqSDIG/"

# Swap F:V between row 78 and row 79
$ws.Range("F78").Value = "Betis"
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = "Valencia"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2.5
$ws.Range("K78").Value = "24/09/2023 17:02"
$ws.Range("L78").Value = 2.07
$ws.Range("M78").Value = "01/10/2023 20:54"
$ws.Range("N78").Value = 3.19
$ws.Range("O78").Value = "24/09/2023 17:02"
$ws.Range("P78").Value = 3.41
$ws.Range("Q78").Value = "01/10/2023 20:54"
$ws.Range("R78").Value = 3.11
$ws.Range("S78").Value = "24/09/2023 17:02"
$ws.Range("T78").Value = 4.04
$ws.Range("U78").Value = "01/10/2023 20:57"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/spain/laliga/betis-valencia/vukArZ2c/"
$ws.Range("F79").Value = "Atl. Madrid"
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = "Cadiz CF"
$ws.Range("I79").Value = 2
$ws.Range("J79").Value = 1.35
$ws.Range("K79").Value = "21/09/2023 22:03"
$ws.Range("L79").Value = 1.34
$ws.Range("M79").Value = "01/10/2023 20:50"
$ws.Range("N79").Value = 5.24
$ws.Range("O79").Value = "21/09/2023 22:03"
$ws.Range("P79").Value = 5.31
$ws.Range("Q79").Value = "01/10/2023 20:59"
$ws.Range("R79").Value = 9.25
$ws.Range("S79").Value = "21/09/2023 22:03"
$ws.Range("T79").Value = 10.48
$ws.Range("U79").Value = "01/10/2023 20:59"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/spain/laliga/atl-madrid-cadiz/E1cOKVAj/"

# Swap F:V between row 179 and row 180
$ws.Range("F179").Value = "Mallorca"
$ws.Range("G179").Value = 3
$ws.Range("H179").Value = "Osasuna"
$ws.Range("I179").Value = 2
$ws.Range("J179").Value = 2.22
$ws.Range("K179").Value = "10/12/2023 10:02"
$ws.Range("L179").Value = 2.26
$ws.Range("M179").Value = "21/12/2023 21:29"
$ws.Range("N179").Value = 3.02
$ws.Range("O179").Value = "10/12/2023 10:02"
$ws.Range("P179").Value = 2.92
$ws.Range("Q179").Value = "21/12/2023 21:29"
$ws.Range("R179").Value = 3.52
$ws.Range("S179").Value = "10/12/2023 10:02"
$ws.Range("T179").Value = 4.11
$ws.Range("U179").Value = "21/12/2023 21:29"
$ws.Range("V179").Value = "https://www.betexplorer.com/football/spain/laliga/mallorca-osasuna/CSRucmzs/"
$ws.Range("F180").Value = "Alaves"
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = "Real Madrid"
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = 4.44
$ws.Range("K180").Value = "10/12/2023 10:02"
$ws.Range("L180").Value = 6.04
$ws.Range("M180").Value = "21/12/2023 21:28"
$ws.Range("N180").Value = 3.85
$ws.Range("O180").Value = "10/12/2023 10:02"
$ws.Range("P180").Value = 4.13
$ws.Range("Q180").Value = "21/12/2023 21:27"
$ws.Range("R180").Value = 1.71
$ws.Range("S180").Value = "10/12/2023 10:02"
$ws.Range("T180").Value = 1.6
$ws.Range("U180").Value = "21/12/2023 21:19"
$ws.Range("V180").Value = "https://www.betexplorer.com/football/spain/laliga/alaves-real-madrid/bqUifoKa/"

# Swap F:V between row 186 and row 187
$ws.Range("F186").Value = "Celta Vigo"
$ws.Range("G186").Value = 2
$ws.Range("H186").Value = "Betis"
$ws.Range("I186").Value = 1
$ws.Range("J186").Value = 1.95
$ws.Range("K186").Value = "17/12/2024 18:03"
$ws.Range("L186").Value = 2.21
$ws.Range("M186").Value = "03/01/2024 19:14"
$ws.Range("N186").Value = 3.47
$ws.Range("O186").Value = "17/12/2024 18:03"
$ws.Range("P186").Value = 3.31
$ws.Range("Q186").Value = "03/01/2024 19:14"
$ws.Range("R186").Value = 3.77
$ws.Range("S186").Value = "17/12/2024 18:03"
$ws.Range("T186").Value = 3.67
$ws.Range("U186").Value = "03/01/2024 19:14"
$ws.Range("V186").Value = "https://www.betexplorer.com/football/spain/laliga/celta-vigo-betis/URKXfIAA/"
$ws.Range("F187").Value = "Real Madrid"
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = "Mallorca"
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = 1.17
$ws.Range("K187").Value = "17/12/2024 18:03"
$ws.Range("L187").Value = 1.22
$ws.Range("M187").Value = "03/01/2024 19:10"
$ws.Range("N187").Value = 6.63
$ws.Range("O187").Value = "17/12/2024 18:03"
$ws.Range("P187").Value = 6.75
$ws.Range("Q187").Value = "03/01/2024 19:14"
$ws.Range("R187").Value = 12.67
$ws.Range("S187").Value = "17/12/2024 18:03"
$ws.Range("T187").Value = 13.93
$ws.Range("U187").Value = "03/01/2024 19:14"
$ws.Range("V187").Value = "https://www.betexplorer.com/football/spain/laliga/real-madrid-mallorca/xhOugduN/"

# Add new rows 189 and 190 (copy formatting from last existing row, then set values)
$ws.Range("A188:V188").Copy()
$ws.Range("A189:V189").PasteSpecial(-4122)
$ws.Range("A188:V188").Copy()
$ws.Range("A190:V190").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 189
$ws.Range("A189").Value = 188
$ws.Range("B189").Value = "spain"
$ws.Range("C189").Value = "laliga"
$ws.Range("D189").Value = "2023-2024"
$ws.Range("E189").Value = 45295.70833333334
$ws.Range("F189").Value = "Osasuna"
$ws.Range("G189").Value = 1
$ws.Range("H189").Value = "Almeria"
$ws.Range("I189").Value = 0
$ws.Range("J189").Value = 1.63
$ws.Range("K189").Value = "17/12/2024 18:03"
$ws.Range("L189").Value = 1.74
$ws.Range("M189").Value = "04/01/2024 16:59"
$ws.Range("N189").Value = 3.92
$ws.Range("O189").Value = "17/12/2024 18:03"
$ws.Range("P189").Value = 3.82
$ws.Range("Q189").Value = "04/01/2024 16:59"
$ws.Range("R189").Value = 4.93
$ws.Range("S189").Value = "17/12/2024 18:03"
$ws.Range("T189").Value = 5.21
$ws.Range("U189").Value = "04/01/2024 16:59"
$ws.Range("V189").Value = "https://www.betexplorer.com/football/spain/laliga/osasuna-almeria/0nTGb0Bo/"

# Row 190
$ws.Range("A190").Value = 189
$ws.Range("B190").Value = "spain"
$ws.Range("C190").Value = "laliga"
$ws.Range("D190").Value = "2023-2024"
$ws.Range("E190").Value = 45295.80208333334
$ws.Range("F190").Value = "Sevilla"
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = "Ath Bilbao"
$ws.Range("I190").Value = 2
$ws.Range("J190").Value = 2.69
$ws.Range("K190").Value = "17/12/2024 18:03"
$ws.Range("L190").Value = 3.09
$ws.Range("M190").Value = "04/01/2024 19:14"
$ws.Range("N190").Value = 3.27
$ws.Range("O190").Value = "17/12/2024 18:03"
$ws.Range("P190").Value = 3.26
$ws.Range("Q190").Value = "04/01/2024 18:50"
$ws.Range("R190").Value = 2.59
$ws.Range("S190").Value = "17/12/2024 18:03"
$ws.Range("T190").Value = 2.52
$ws.Range("U190").Value = "04/01/2024 19:14"
$ws.Range("V190").Value = "https://www.betexplorer.com/football/spain/laliga/sevilla-ath-bilbao/vwSKcKQi/"
